$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "na?o informado" row (row 337) entirely; Excel will shift
# every subsequent row up by one, which also reduces the used range from
# A1:C586 down to A1:C585.
$ws.Rows.Item(337).Delete()
